$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 2 ("Features: ...") : remove "humidity, " from the
#     middle of the line, splitting the remaining text into three runs
#     ("Features: N, P, K, " / "temperature, " / "pH, rainfall") the
#     same way PowerPoint leaves the run boundaries behind after an
#     in-place delete. ---
$para2 = $tr.Paragraphs(2, 1)
$full = $para2.Text
$cut = $full.IndexOf("humidity, ")
$delStart = $para2.Start + $cut
$toDelete = $tr.Characters($delStart, "humidity, ".Length)
$toDelete.Text = ""

# Re-assigning the (unchanged) text of the middle segment back onto
# itself forces PowerPoint's run-splitting without touching formatting,
# producing the extra run boundary before "temperature, ".
$para2 = $tr.Paragraphs(2, 1)
$segA = "Features: N, P, K, "
$midStart = $para2.Start + $segA.Length
$midLen = "temperature, ".Length
$midRange = $tr.Characters($midStart, $midLen)
$midRange.Text = $midRange.Text

# --- Shape resize: the text box was manually resized (both width and
#     height shrank) after the edit. ---
$sh.Width = 357.305591
$sh.Height = 72.70316
